$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new keyword ("Architectuur") for Hans Krusemann (row 14), marking the row as seen ---
# Copy the existing "seen" checkmark style from A10 onto A14 (column A = "Gezien?")
$ws.Range("A10").Copy($ws.Range("A14"))

# Copy the keyword-cell formatting (bold, centered, rotated) from E23 onto F14 ("Architectuur" column),
# then fill in the actual keyword text.
$ws.Range("E23").Copy($ws.Range("F14"))
$ws.Range("F14").Value = "Architectuur"

# Row grew slightly taller to fit the new content
$ws.Rows.Item(14).RowHeight = 58

# --- Mark two more rows as reviewed (Anke Spijker / row 8, Marika Beckers-van Hout / row 18) ---
$ws.Range("A8").Value = "."
$ws.Range("A18").Value = "."

# --- Totals row: exclude the blank separator row from the "Gezien?" ratio ---
$ws.Range("A31").Formula = "=COUNTA(A7:A29)/COUNTA(B7:B29)"

# --- Update the last worked-on cell/selection ---
$ws.Range("Y21").Select()
